{"js": "// Replace the three-digit x one-digit multiplication problems in the table\n// with the new values, per the commit's diff. Every \"old\" text below is\n// unique within the document, so a literal (non-wildcard) search safely\n// targets exactly one run each.\nconst replacements = [\n  [\"382\u00d79=\", \"286\u00d72=\"],\n  [\"452\u00d74=\", \"357\u00d77=\"],\n  [\"529\u00d76=\", \"103\u00d75=\"],\n  [\"524\u00d79=\", \"734\u00d74=\"],\n  [\"947\u00d76=\", \"518\u00d77=\"],\n  [\"437\u00d74=\", \"551\u00d77=\"],\n  [\"909\u00d77=\", \"686\u00d78=\"],\n  [\"908\u00d77=\", \"969\u00d79=\"],\n  [\"736\u00d75=\", \"381\u00d79=\"],\n  [\"158\u00d74=\", \"481\u00d75=\"],\n  [\"617\u00d72=\", \"717\u00d77=\"],\n  [\"993\u00d78=\", \"405\u00d72=\"],\n  [\"152\u00d79=\", \"850\u00d73=\"],\n  [\"664\u00d75=\", \"587\u00d73=\"],\n  [\"551\u00d79=\", \"249\u00d75=\"],\n  [\"221\u00d72=\", \"379\u00d75=\"],\n  [\"270\u00d72=\", \"161\u00d76=\"],\n  [\"329\u00d74=\", \"534\u00d73=\"],\n  [\"307\u00d74=\", \"130\u00d78=\"],\n  [\"683\u00d79=\", \"992\u00d72=\"],\n  [\"838\u00d76=\", \"881\u00d78=\"],\n  [\"606\u00d75=\", \"548\u00d76=\"],\n  [\"735\u00d78=\", \"674\u00d74=\"],\n  [\"611\u00d79=\", \"341\u00d73=\"],\n  [\"612\u00d79=\", \"514\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems in the table\n# with the new values, per the commit's diff. Every \"old\" text below is\n# unique within the document, so Find/Replace safely targets exactly one\n# run each; wdReplaceAll (2) is used defensively but only ever matches once.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"382\u00d79=\", \"286\u00d72=\"),\n    @(\"452\u00d74=\", \"357\u00d77=\"),\n    @(\"529\u00d76=\", \"103\u00d75=\"),\n    @(\"524\u00d79=\", \"734\u00d74=\"),\n    @(\"947\u00d76=\", \"518\u00d77=\"),\n    @(\"437\u00d74=\", \"551\u00d77=\"),\n    @(\"909\u00d77=\", \"686\u00d78=\"),\n    @(\"908\u00d77=\", \"969\u00d79=\"),\n    @(\"736\u00d75=\", \"381\u00d79=\"),\n    @(\"158\u00d74=\", \"481\u00d75=\"),\n    @(\"617\u00d72=\", \"717\u00d77=\"),\n    @(\"993\u00d78=\", \"405\u00d72=\"),\n    @(\"152\u00d79=\", \"850\u00d73=\"),\n    @(\"664\u00d75=\", \"587\u00d73=\"),\n    @(\"551\u00d79=\", \"249\u00d75=\"),\n    @(\"221\u00d72=\", \"379\u00d75=\"),\n    @(\"270\u00d72=\", \"161\u00d76=\"),\n    @(\"329\u00d74=\", \"534\u00d73=\"),\n    @(\"307\u00d74=\", \"130\u00d78=\"),\n    @(\"683\u00d79=\", \"992\u00d72=\"),\n    @(\"838\u00d76=\", \"881\u00d78=\"),\n    @(\"606\u00d75=\", \"548\u00d76=\"),\n    @(\"735\u00d78=\", \"674\u00d74=\"),\n    @(\"611\u00d79=\", \"341\u00d73=\"),\n    @(\"612\u00d79=\", \"514\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
